$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-15 from 45243 to 45244 (2023-11-13 -> 2023-11-14)
for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value = 45244
    }
}
